# Update the two "Data" sheet notes cells whose text changed in this revision,
# then leave the selection on B5 (matching the saved cursor position).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# D3: clarify that the mobility "countries/regions, sub-regions, and cities"
# note also covers county-level data.
$ws.Range("D3").Value = "countries/regions, sub-regions, and cities, county"

# B5: fix the cases/deaths/recoveries formula note (was "deaths - cases (?)").
$ws.Range("B5").Value = "cases, deaths, recoveries = cases - deaths"

$ws.Range("B5").Select()
